$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column F
$ws.Range("F1").Value = "time_taken"
# Copy just the formatting from the existing header style (E1) onto F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for column F (unstyled, like the other data cells)
$ws.Range("F2").Value = "2021-10-05 10:52:34.924092"
$ws.Range("F3").Value = "2021-10-05 10:52:34.924105"
$ws.Range("F4").Value = "2021-10-05 10:52:34.924108"
$ws.Range("F5").Value = "2021-10-05 10:52:34.924112"
